$d = $word.ActiveDocument

$replacements = @(
    @("2024-10-30 Wednesday", "2024-10-31 Thursday"),
    @("63×16=1008", "52×19=988"),
    @("43×30=1290", "33×91=3003"),
    @("51×26=1326", "86×62=5332"),
    @("87×36=3132", "57×59=3363"),
    @("50×33=1650", "50×85=4250"),
    @("80×15=1200", "17×30=510"),
    @("46×97=4462", "12×88=1056"),
    @("61×96=5856", "93×13=1209"),
    @("34×61=2074", "15×39=585"),
    @("79×47=3713", "20×26=520"),
    @("37×40=1480", "31×98=3038"),
    @("17×16=272", "54×56=3024"),
    @("86×97=8342", "93×35=3255"),
    @("90×55=4950", "80×11=880"),
    @("93×14=1302", "84×34=2856"),
    @("33×55=1815", "27×23=621"),
    @("58×11=638", "47×98=4606"),
    @("35×99=3465", "15×95=1425"),
    @("96×22=2112", "62×53=3286"),
    @("26×21=546", "63×59=3717"),
    @("22×61=1342", "37×59=2183"),
    @("23×29=667", "26×77=2002"),
    @("54×31=1674", "41×40=1640"),
    @("13×89=1157", "25×73=1825"),
    @("90×44=3960", "55×77=4235")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
